$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric are preserved as text (matches source formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.977.06'
$ws.Range("E2").Value = '  -3.93%  '
$ws.Range("D3").Value = '3.515.95'
$ws.Range("E3").Value = '  -4.64%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '579.99'
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("D6").Value = '174.34'
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.510.36'
$ws.Range("E8").Value = '  -4.65%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D11").Value = '6.71'
$ws.Range("E11").Value = '  +8.44%  '
$ws.Range("D12").Value = '0.599'
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").Value = '47.27'
$ws.Range("E13").Value = '  -5.33%  '
$ws.Range("E14").Value = '  -3.33%  '
$ws.Range("D15").Value = '670.98'
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("D16").Value = '4.080.53'
$ws.Range("E16").Value = '  -4.68%  '
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").Value = '3.512.36'
$ws.Range("E18").Value = '  -4.67%  '
$ws.Range("D19").Value = '68.954.76'
$ws.Range("E19").Value = '  -4.16%  '
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("D21").Value = '17.55'
$ws.Range("E21").Value = '  -2.71%  '
$ws.Range("D22").Value = '11.23'
$ws.Range("E22").Value = '  -3.53%  '
$ws.Range("D23").Value = '0.907'
$ws.Range("E23").Value = '  -3.49%  '
$ws.Range("D24").Value = '16.30'
$ws.Range("E24").Value = '  -8.43%  '
$ws.Range("D25").Value = '98.32'
$ws.Range("E25").Value = '  -4.98%  '
$ws.Range("D26").Value = '3.86'
$ws.Range("E26").Value = '  -4.20%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '2.66'
$ws.Range("E28").Value = '  -6.46%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '9.45'
$ws.Range("E29").Value = '  -6.98%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '32.95'
$ws.Range("E30").Value = '  -7.06%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '8.76'
$ws.Range("E31").Value = '  -4.71%  '
$ws.Range("E32").Value = '  -7.52%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '7.32'
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = '1.36'
$ws.Range("E34").Value = '  -4.68%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").Value = '577.71'
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("B36").Value = 'Cosmos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D36").Value = '10.94'
$ws.Range("E36").Value = '  -3.24%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = '3.59'
$ws.Range("E37").Value = '  -14.34%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.105'
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '57.13'
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '0.338'
$ws.Range("E41").Value = '  -3.30%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0438'
$ws.Range("E42").Value = '  -5.23%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = '0.137'
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.419.75'
$ws.Range("E44").Value = '  -9.01%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '33.48'
$ws.Range("E45").Value = '  -5.45%  '
$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").Value = '0.0₃0706'
$ws.Range("E46").Value = '  -8.90%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").Value = '2.93'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").Value = '2.60'
$ws.Range("E48").Value = '  -6.87%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '131.18'
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.150'
$ws.Range("E51").Value = '  -0.29%  '
